$d = $word.ActiveDocument

# --- Step 1: Fix the trailing-space run at the end of the last paragraph ---
# ("go down that route ... buttons. " -> "go down that route ... buttons." , drop xml:space=preserve)
$lastPara = $d.Paragraphs.Last
$paraRange = $lastPara.Range
$paraText = $paraRange.Text
$searchText = "go down that route and need it if I can get the radio buttons to work as well as separating the Prev and next buttons."
$idx = $paraText.IndexOf($searchText)
if ($idx -lt 0) {
    throw "Could not locate target run text to fix"
}
$runStart = $paraRange.Start + $idx
$runEnd = $paraRange.End - 1
$runRange = $d.Range($runStart, $runEnd)

$fixXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>go down that route and need it if I can get the radio buttons to work as well as separating the Prev and next buttons.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$runRange.InsertXML($fixXml)

# --- Step 2: Append the two new paragraphs at the very end of the document ---
$endPos = $d.Content.End
$insertPoint = $d.Range($endPos, $endPos)

$newParasXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Captain''s log. Stardate 02.</w:t></w:r><w:r><w:t>25</w:t></w:r><w:r><w:t>.24</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>F</w:t></w:r><w:r><w:t>ix</w:t></w:r><w:r><w:t>ed</w:t></w:r><w:r><w:t xml:space="preserve"> the cards </w:t></w:r><w:r><w:t>so</w:t></w:r><w:r><w:t xml:space="preserve"> that the</w:t></w:r><w:r><w:t>y</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">don’t break the code and cycle though the questions. </w:t></w:r><w:r><w:t xml:space="preserve">Added </w:t></w:r><w:r><w:t>radio</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>buttons</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t xml:space="preserve">was able to </w:t></w:r><w:r><w:t>make it so the user can select an answer and it will show if it is wrong or right. Fixed the issue where the fill in the blank wasn’t showing</w:t></w:r><w:r><w:t xml:space="preserve">, got the first set to show. Then fixed that so now they all show up under the fill in the blank button. Fixed the issue that the answers </w:t></w:r><w:r><w:t>were</w:t></w:r><w:r><w:t xml:space="preserve"> still showing when </w:t></w:r><w:r><w:t xml:space="preserve">you hit the next card or a different level. </w:t></w:r><w:r><w:t>Forgot that</w:t></w:r><w:r><w:t xml:space="preserve"> I needed to also add that to the PREV </w:t></w:r><w:r><w:t>button,</w:t></w:r><w:r><w:t xml:space="preserve"> so I fixed that </w:t></w:r><w:r><w:t xml:space="preserve">quickly as well. Broke my code multiple times when trying to fix it but in the end worked line by line and was able to successfully fix it. </w:t></w:r><w:r><w:t>Created props so that my variables aren’t being re-declared</w:t></w:r><w:r><w:t xml:space="preserve"> fixing an issue</w:t></w:r><w:r><w:t xml:space="preserve"> and cleaned some redundant code. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertPoint.InsertXML($newParasXml)
